$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1727.1052
$ws.Range("I40").Value = 1736
$ws.Range("J40").Value = 1710
$ws.Range("K40").Value = 1736
$ws.Range("L40").Value = 1710
$ws.Range("M40").Value = -1561
$ws.Range("N40").Value = -2060

$ws.Range("H43").Value = 752.75
$ws.Range("I43").Value = 668.7143
$ws.Range("J43").Value = 870.4
$ws.Range("K43").Value = 668.7143
$ws.Range("L43").Value = 870.4
$ws.Range("M43").Value = -599.7143
$ws.Range("N43").Value = -1008.4

$ws.Range("H48").Value = 3380.353
$ws.Range("I48").Value = 1704
$ws.Range("J48").Value = 3896.1538
$ws.Range("K48").Value = 5112
$ws.Range("L48").Value = 11688.4614
$ws.Range("M48").Value = -4820
$ws.Range("N48").Value = -12272.4614

$ws.Range("H56").Value = 3380.353
$ws.Range("I56").Value = 1704
$ws.Range("J56").Value = 3896.1538
$ws.Range("K56").Value = 5112
$ws.Range("L56").Value = 11688.4614
$ws.Range("M56").Value = -4578
$ws.Range("N56").Value = -12756.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 325.4074
$ws.Range("I2").Value = 325.4074
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 325.4074
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -212.4074
$ws.Range("N2").ClearContents()

$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -384
$ws.Range("N4").ClearContents()

$ws.Range("H32").Value = 15155662
$ws.Range("I32").Value = 18183652
$ws.Range("J32").Value = 15717.272
$ws.Range("K32").Value = 18183652
$ws.Range("L32").Value = 15717.272
$ws.Range("M32").Value = -18183365
$ws.Range("N32").Value = -16291.272

$ws.Range("H61").Value = 2030.9143
$ws.Range("I61").Value = 2400.5
$ws.Range("J61").Value = 1538.1333
$ws.Range("K61").Value = 2400.5
$ws.Range("L61").Value = 1538.1333
$ws.Range("M61").Value = -2188.5
$ws.Range("N61").Value = -1962.1333

$ws.Range("H102").Value = 2769.9
$ws.Range("I102").Value = 1926.8572
$ws.Range("J102").Value = 4737
$ws.Range("K102").Value = 1926.8572
$ws.Range("L102").Value = 4737
$ws.Range("M102").Value = -304.8571999999999
$ws.Range("N102").Value = -7981

$ws.Range("H116").Value = 325.4074
$ws.Range("I116").Value = 325.4074
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 325.4074
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1968.5926
$ws.Range("N116").ClearContents()

$ws.Range("H136").Value = 2030.9143
$ws.Range("I136").Value = 2400.5
$ws.Range("J136").Value = 1538.1333
$ws.Range("K136").Value = 7201.5
$ws.Range("L136").Value = 4614.3999
$ws.Range("M136").Value = -4651.5
$ws.Range("N136").Value = -9714.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 325.4074
$ws.Range("I3").Value = 325.4074
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 325.4074
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -211.4074
$ws.Range("N3").ClearContents()

$ws.Range("H94").Value = 796.73914
$ws.Range("I94").Value = 779.5833
$ws.Range("J94").Value = 815.4545000000001
$ws.Range("K94").Value = 779.5833
$ws.Range("L94").Value = 815.4545000000001
$ws.Range("M94").Value = -328.5833
$ws.Range("N94").Value = -1717.4545

$ws.Range("H107").Value = 1330.875
$ws.Range("I107").Value = 1019.6
$ws.Range("J107").Value = 6000
$ws.Range("K107").Value = 1019.6
$ws.Range("L107").Value = 6000
$ws.Range("M107").Value = 900.4
$ws.Range("N107").Value = -9840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1105.127
$ws.Range("I31").Value = 1263.826
$ws.Range("J31").Value = 1013.875
$ws.Range("K31").Value = 1263.826
$ws.Range("L31").Value = 1013.875
$ws.Range("M31").Value = -968.826
$ws.Range("N31").Value = -1603.875

$ws.Range("H34").Value = 1105.127
$ws.Range("I34").Value = 1263.826
$ws.Range("J34").Value = 1013.875
$ws.Range("K34").Value = 1263.826
$ws.Range("L34").Value = 1013.875
$ws.Range("M34").Value = -1061.826
$ws.Range("N34").Value = -1417.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 729.4286
$ws.Range("I4").Value = 356.77777
$ws.Range("J4").Value = 1400.2
$ws.Range("K4").Value = 1070.33331
$ws.Range("L4").Value = 4200.6
$ws.Range("M4").Value = -958.33331
$ws.Range("N4").Value = -4424.6

$ws.Range("H64").Value = 2062.5625
$ws.Range("I64").Value = 700.1667
$ws.Range("J64").Value = 2880
$ws.Range("K64").Value = 2100.5001
$ws.Range("L64").Value = 8640
$ws.Range("M64").Value = -1830.5001
$ws.Range("N64").Value = -9180

$ws.Range("H67").Value = 2062.5625
$ws.Range("I67").Value = 700.1667
$ws.Range("J67").Value = 2880
$ws.Range("K67").Value = 2100.5001
$ws.Range("L67").Value = 8640
$ws.Range("M67").Value = -1164.5001
$ws.Range("N67").Value = -10512

$ws.Range("H80").Value = 862.4286
$ws.Range("I80").Value = 777
$ws.Range("J80").Value = 876.6667
$ws.Range("K80").Value = 2331
$ws.Range("L80").Value = 2630.0001
$ws.Range("M80").Value = -1395
$ws.Range("N80").Value = -4502.0001

$ws.Range("H83").Value = 862.4286
$ws.Range("I83").Value = 777
$ws.Range("J83").Value = 876.6667
$ws.Range("K83").Value = 6993
$ws.Range("L83").Value = 7890.0003
$ws.Range("M83").Value = -2313
$ws.Range("N83").Value = -17250.0003

$ws.Range("H114").Value = 1782.5
$ws.Range("I114").Value = 1334.5454
$ws.Range("J114").Value = 2486.4285
$ws.Range("K114").Value = 4003.6362
$ws.Range("L114").Value = 7459.2855
$ws.Range("M114").Value = -749.6361999999999
$ws.Range("N114").Value = -13967.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 25022.857
$ws.Range("J93").Value = 25022.857
$ws.Range("L93").Value = 25022.857
$ws.Range("N93").Value = -28766.857

$ws.Range("H102").Value = 2346.8667
$ws.Range("I102").Value = 2346.3845
$ws.Range("J102").Value = 2350
$ws.Range("K102").Value = 2346.3845
$ws.Range("L102").Value = 2350
$ws.Range("M102").Value = -724.3845000000001
$ws.Range("N102").Value = -5594

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 680.64703
$ws.Range("I55").Value = 548.875
$ws.Range("J55").Value = 797.7778
$ws.Range("K55").Value = 548.875
$ws.Range("L55").Value = 797.7778
$ws.Range("M55").Value = -375.875
$ws.Range("N55").Value = -1143.7778

$ws.Range("H136").Value = 1320.1818
$ws.Range("I136").Value = 1236.3846
$ws.Range("J136").Value = 1631.4286
$ws.Range("K136").Value = 3709.1538
$ws.Range("L136").Value = 4894.2858
$ws.Range("M136").Value = -1159.1538
$ws.Range("N136").Value = -9994.2858
